# "integrate database all in one"
# A new reporting period (12 ماهه منتهی به 1399/05 / published 1400-02-03 (7))
# is inserted into column E, the old column D figures are replaced with what
# used to be in column E, and the brand-new period's figures land in column E.
# Also drops the now-unused "-" placeholder string from D15 (replaced with a
# real numeric value), which is what it used to be in column E.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# --- Row 8: financial period headers ---
$ws.Range("D8").Value = $ws.Range("E8").Value()
$ws.Range("E8").Value = "12 ماهه منتهی به 1399/05"

# --- Row 9: publish dates ---
$ws.Range("D9").Value = $ws.Range("E9").Value()
$ws.Range("E9").Value = "1400-02-03 (7)"
$ws.Range("F9").Value = "1400-10-08 (8)"
$ws.Range("H9").Value = "1402-03-30 (5)"

# --- Row 11: فروش (Sales) ---
$ws.Range("D11").Value = 1709782
$ws.Range("E11").Value = 1981727

# --- Row 12: بهای تمام شده کالای فروش رفته ---
$ws.Range("D12").Value = -875747
$ws.Range("E12").Value = -1036745

# --- Row 13: سود (زیان) ناخالص ---
$ws.Range("D13").Value = 834035
$ws.Range("E13").Value = 944982

# --- Row 14: هزینه های عمومی, اداری و تشکیلاتی ---
$ws.Range("D14").Value = -104958
$ws.Range("E14").Value = -144083

# --- Row 15: هزینه کاهش ارزش دریافتنی‌‏ها (هزینه استثنایی) ---
$ws.Range("D15").Value = -6877
$ws.Range("E15").Value = 0

# --- Row 16: خالص سایر درامدها (هزینه ها) ی عملیاتی ---
$ws.Range("D16").Value = 8090
$ws.Range("E16").Value = 50419

# --- Row 17: سود (زیان) عملیاتی ---
$ws.Range("D17").Value = 730290
$ws.Range("E17").Value = 851318

# --- Row 18: هزینه های مالی ---
$ws.Range("D18").Value = -71889
$ws.Range("E18").Value = -81222

# --- Row 19: خالص سایر درامدها و هزینه های غیرعملیاتی ---
$ws.Range("D19").Value = 42612
$ws.Range("E19").Value = 200918

# --- Row 20: سود (زیان) خالص عملیات در حال تداوم قبل از مالیات ---
$ws.Range("D20").Value = 701013
$ws.Range("E20").Value = 971014

# --- Row 21: مالیات ---
$ws.Range("D21").Value = -74023
$ws.Range("E21").Value = -70616

# --- Row 22: سود (زیان) خالص عملیات در حال تداوم ---
$ws.Range("D22").Value = 626990
$ws.Range("E22").Value = 900398

# --- Row 23: سود (زیان) عملیات متوقف شده پس از اثر مالیاتی --- (unchanged 0/0)
$ws.Range("D23").Value = 0
$ws.Range("E23").Value = 0

# --- Row 24: سود (زیان) خالص ---
$ws.Range("D24").Value = 626990
$ws.Range("E24").Value = 900398

# --- Row 25: سود هر سهم پس از کسر مالیات ---
$ws.Range("D25").Value = 1254
$ws.Range("E25").Value = 1801

# --- Row 26: سرمایه ---
$ws.Range("D26").Value = 500000
$ws.Range("E26").Value = 500000

# --- Row 27: سود هر سهم بر اساس آخرین سرمایه ---
$ws.Range("D27").Value = 627
$ws.Range("E27").Value = 900
